$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy row-121 formatting down onto the two new rows (122, 123) ---
# Column A carries the bold/centered/bordered "index" style.
$ws.Range("A121").Copy()
$ws.Range("A122:A123").PasteSpecial(-4122)

# Column E carries the date/time number format.
$ws.Range("E121").Copy()
$ws.Range("E122:E123").PasteSpecial(-4122)

# Column D ("temporada") is stored as text ("2023"); pasting the value from
# the row above (same literal text) keeps it text instead of Excel's
# auto-number coercion that a plain .Value assignment would trigger.
# (done one cell at a time - a multi-cell PasteSpecial destination only
# fills the first cell of the range in this host)
$ws.Range("D121").Copy()
$ws.Range("D122").PasteSpecial(-4163)
$ws.Range("D121").Copy()
$ws.Range("D123").PasteSpecial(-4163)

$excel.CutCopyMode = $false

# --- Row 122 : Guairena 4 x 1 Resistencia ---
$ws.Range("A122").Value = 121
$ws.Range("B122").Value = "paraguay"
$ws.Range("C122").Value = "primera-division"
$ws.Range("E122").Value = 45253.97916666666
$ws.Range("F122").Value = "Guairena"
$ws.Range("G122").Value = 4
$ws.Range("H122").Value = "Resistencia"
$ws.Range("I122").Value = 1
$ws.Range("J122").Value = 1.89
$ws.Range("K122").Value = "19/11/2023 01:12"
$ws.Range("L122").Value = 1.56
$ws.Range("M122").Value = "23/11/2023 23:22"
$ws.Range("N122").Value = 3.6
$ws.Range("O122").Value = "19/11/2023 01:12"
$ws.Range("P122").Value = 4.31
$ws.Range("Q122").Value = "23/11/2023 23:22"
$ws.Range("R122").Value = 4.19
$ws.Range("S122").Value = "19/11/2023 01:12"
$ws.Range("T122").Value = 6.14
$ws.Range("U122").Value = "23/11/2023 23:22"
$ws.Range("V122").Value = "https://www.betexplorer.com/football/paraguay/primera-division/guairena-fc-resistencia/MZ0rrU1l/"

# --- Row 123 : Tacuary 1 x 1 Sp. Luqueno ---
$ws.Range("A123").Value = 122
$ws.Range("B123").Value = "paraguay"
$ws.Range("C123").Value = "primera-division"
$ws.Range("E123").Value = 45253.97916666666
$ws.Range("F123").Value = "Tacuary"
$ws.Range("G123").Value = 1
$ws.Range("H123").Value = "Sp. Luqueno"
$ws.Range("I123").Value = 1
$ws.Range("J123").Value = 3.82
$ws.Range("K123").Value = "19/11/2023 01:12"
$ws.Range("L123").Value = 3.76
$ws.Range("M123").Value = "23/11/2023 23:29"
$ws.Range("N123").Value = 3.57
$ws.Range("O123").Value = "19/11/2023 01:12"
$ws.Range("P123").Value = 3.33
$ws.Range("Q123").Value = "23/11/2023 23:29"
$ws.Range("R123").Value = 1.92
$ws.Range("S123").Value = "19/11/2023 01:12"
$ws.Range("T123").Value = 2.14
$ws.Range("U123").Value = "23/11/2023 23:29"
$ws.Range("V123").Value = "https://www.betexplorer.com/football/paraguay/primera-division/tacuary-sp-luqueno/rq2vqAnr/"

Write-Host "Rows 122-123 added"
